# Generate Report for Handoff
# Refresh the "Priority" and "Latest Handoff Datetime" columns for the rows
# that had just been handed off (previously "low" priority / stale handoff
# timestamp) in both the zh-cn and de-de localization-status sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-10-18 13:44:43"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-10-18 13:44:54"

# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# de-de handoff datetime for these rows, so refresh it too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4:G7").Value = "2016-10-18 13:44:54"
